$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 11
$ws1.Range("F16").Value = 4371
$ws1.Range("F22").Value = 994
$ws1.Range("F23").Value = 1819
$ws1.Range("F34").Value = 27

# Sheet "全部类型" (All types) - update "想去人数" (want-to-go count) column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 11
$ws4.Range("F17").Value = 4371
$ws4.Range("F23").Value = 994
$ws4.Range("F24").Value = 1819
$ws4.Range("F35").Value = 27
